# Apply updates described by the commit "Update gh-pages to output generated at 456a3b4"
# Workbook has 4 sheets: 展览 (Exhibition), 演出 (Performance), 本地生活 (Local life), 全部类型 (All types)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" - bump the "想去人数" (want-to-go count) column F for many rows
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 288
$ws.Cells.Item(5, 6).Value = 2927
$ws.Cells.Item(6, 6).Value = 19
$ws.Cells.Item(7, 6).Value = 237
$ws.Cells.Item(10, 6).Value = 6848
$ws.Cells.Item(11, 6).Value = 37
$ws.Cells.Item(12, 6).Value = 55
$ws.Cells.Item(13, 6).Value = 346
$ws.Cells.Item(14, 6).Value = 599
$ws.Cells.Item(15, 6).Value = 1483
$ws.Cells.Item(16, 6).Value = 1109
$ws.Cells.Item(18, 6).Value = 1468
$ws.Cells.Item(20, 6).Value = 104
$ws.Cells.Item(21, 6).Value = 1106
$ws.Cells.Item(22, 6).Value = 119
$ws.Cells.Item(23, 6).Value = 172
$ws.Cells.Item(24, 6).Value = 333
$ws.Cells.Item(25, 6).Value = 1692
$ws.Cells.Item(26, 6).Value = 1682
$ws.Cells.Item(28, 6).Value = 1029
$ws.Cells.Item(29, 6).Value = 34
$ws.Cells.Item(30, 6).Value = 1658
$ws.Cells.Item(31, 6).Value = 1204
$ws.Cells.Item(32, 6).Value = 135
$ws.Cells.Item(34, 6).Value = 24
$ws.Cells.Item(36, 6).Value = 417
$ws.Cells.Item(37, 6).Value = 6
$ws.Cells.Item(38, 6).Value = 2447
$ws.Cells.Item(39, 6).Value = 2706
$ws.Cells.Item(40, 6).Value = 70
$ws.Cells.Item(44, 6).Value = 22
$ws.Cells.Item(46, 6).Value = 119
$ws.Cells.Item(48, 6).Value = 151
$ws.Cells.Item(49, 6).Value = 412

# ---------------------------------------------------------------------------
# Sheet "演出" - bump column F counts, and mark row 9 as "不可售" (unavailable)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(8, 6).Value = 211
$ws.Cells.Item(9, 7).Value = "不可售"
$ws.Cells.Item(12, 6).Value = 179
$ws.Cells.Item(17, 6).Value = 163
$ws.Cells.Item(20, 6).Value = 44
$ws.Cells.Item(23, 6).Value = 464

# ---------------------------------------------------------------------------
# Sheet "本地生活" - the "全职高手" entry (row 6) was removed entirely; every
# following row shifts up by one. Deleting the row handles the shift of all
# columns (and preserves per-cell styling). Afterwards the sequential index
# in column A is restored, and a handful of F-column counts get bumped to
# match newer scrape values.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Rows.Item(6).Delete()

# Restore sequential index column (A) which is independent of row content
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(14, 1).Value = 13

# Updated "想去人数" counts for the rows that shifted up
$ws.Cells.Item(8, 6).Value = 2713
$ws.Cells.Item(9, 6).Value = 998
$ws.Cells.Item(10, 6).Value = 906
$ws.Cells.Item(12, 6).Value = 251
$ws.Cells.Item(13, 6).Value = 1412
$ws.Cells.Item(14, 6).Value = 7294

# ---------------------------------------------------------------------------
# Sheet "全部类型" - bump column F counts (same underlying rows as other
# sheets, combined into a single "all types" listing)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(4, 6).Value = 288
$ws.Cells.Item(5, 6).Value = 2927
$ws.Cells.Item(6, 6).Value = 237
$ws.Cells.Item(9, 6).Value = 2713
$ws.Cells.Item(10, 6).Value = 6848
$ws.Cells.Item(11, 6).Value = 998
$ws.Cells.Item(12, 6).Value = 906
$ws.Cells.Item(13, 6).Value = 37
$ws.Cells.Item(14, 6).Value = 346
$ws.Cells.Item(16, 6).Value = 251
$ws.Cells.Item(17, 6).Value = 1412
$ws.Cells.Item(18, 6).Value = 599
$ws.Cells.Item(20, 6).Value = 1468
$ws.Cells.Item(22, 6).Value = 104
$ws.Cells.Item(23, 6).Value = 1106
$ws.Cells.Item(24, 6).Value = 119
$ws.Cells.Item(25, 6).Value = 333
$ws.Cells.Item(27, 6).Value = 1692
$ws.Cells.Item(28, 6).Value = 1029
$ws.Cells.Item(30, 6).Value = 34
$ws.Cells.Item(31, 6).Value = 1658
$ws.Cells.Item(32, 6).Value = 1204
$ws.Cells.Item(33, 6).Value = 135
$ws.Cells.Item(34, 6).Value = 24
$ws.Cells.Item(36, 6).Value = 464
$ws.Cells.Item(37, 6).Value = 417
$ws.Cells.Item(39, 6).Value = 2447
$ws.Cells.Item(40, 6).Value = 2706
$ws.Cells.Item(41, 6).Value = 70
$ws.Cells.Item(45, 6).Value = 119
$ws.Cells.Item(48, 6).Value = 412
